$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetimes for the
# 9aa91356-fb3d-41bc-9c59-cdff5b18e997 row group (rows 4 and 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E5").Value = "2016-03-17 18:16:14"
$wsZhCn.Range("H4:H5").Value = "2016-03-17 18:16:45"

# de-de sheet: same update, different timestamps
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E5").Value = "2016-03-17 18:16:22"
$wsDeDe.Range("H4:H5").Value = "2016-03-17 18:16:51"
